# doc/perf.xlsx: correct classpath for JOPtimizer (BCEL 5.2 class loading bug)
# The "JOP" column (F) in the "compare" sheet holds stale numbers; update
# them to the corrected measurements. Downstream ratio formulas (row 7-9)
# and the dependent charts recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("compare")

$ws.Range("F2").Value = 24058
$ws.Range("F3").Value = 10144
$ws.Range("F4").Value = 24308

# Select/activate the "compare" sheet (it becomes the active tab, matching
# the workbook's bookViews no longer pinning "trend" as active) and land
# the selection on the corrected cell.
$ws.Activate()
$ws.Range("F4").Select()
